$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.575056
$ws.Range("H2").Value = 5.150112
$ws.Range("I2").Value = 0.2861985448204127
$ws.Range("J2").Value = 0.2146404527606434
$ws.Range("Q2").Value = 0.809676574784
$ws.Range("R2").Value = 4.858059448704
$ws.Range("S2").Value = 0.2861985448204127
$ws.Range("T2").Value = 0.2146404527606434

# Row 3
$ws.Range("I3").Value = 0.005499560551228268
$ws.Range("J3").Value = 0.006186761889936853
$ws.Range("S3").Value = 0.005499560551228268
$ws.Range("T3").Value = 0.006186761889936853

# Row 4
$ws.Range("G4").Value = 1.021765666666667
$ws.Range("H4").Value = 3.065297
$ws.Range("I4").Value = 0.1135617427145114
$ws.Range("J4").Value = 0.1277519277106676
$ws.Range("Q4").Value = 0.3212744597471111
$ws.Range("R4").Value = 2.891470137724
$ws.Range("S4").Value = 0.1135617427145114
$ws.Range("T4").Value = 0.1277519277106676

# Row 5
$ws.Range("G5").Value = 0.4231475
$ws.Range("H5").Value = 0.846295
$ws.Range("I5").Value = 0.04702973401137512
$ws.Range("J5").Value = 0.0352709109955412
$ws.Range("Q5").Value = 0.1330505505233333
$ws.Range("R5").Value = 0.7983033031400001
$ws.Range("S5").Value = 0.04702973401137512
$ws.Range("T5").Value = 0.0352709109955412

# Row 6
$ws.Range("G6").Value = 3.063320333333333
$ws.Range("H6").Value = 9.189961
$ws.Range("I6").Value = 0.3404655361742741
$ws.Range("J6").Value = 0.3830086394029208
$ws.Range("Q6").Value = 0.9632018546235557
$ws.Range("R6").Value = 8.668816691612001
$ws.Range("S6").Value = 0.3404655361742741
$ws.Range("T6").Value = 0.3830086394029208

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 1.864674666666666
$ws.Range("H7").Value = 5.594023999999999
$ws.Range("I7").Value = 0.2072448817281985
$ws.Range("J7").Value = 0.2331413072402901
$ws.Range("Q7").Value = 0.5863108985564445
$ws.Range("R7").Value = 5.276798087007999
$ws.Range("S7").Value = 0.2072448817281985
$ws.Range("T7").Value = 0.2331413072402901
